# Auto-generated Excel COM-interop script applying the diff
# Updates cached numeric values across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H33").Value = 545.9091
$ws.Range("I33").Value = 324.11765
$ws.Range("K33").Value = 324.11765
$ws.Range("M33").Value = -95.11765000000003
$ws.Range("H55").Value = 381.16666
$ws.Range("J55").Value = 499.25
$ws.Range("L55").Value = 499.25
$ws.Range("N55").Value = -927.25
$ws.Range("H70").Value = 1635.6
$ws.Range("I70").Value = 1087.5
$ws.Range("J70").Value = 2001
$ws.Range("K70").Value = 3262.5
$ws.Range("L70").Value = 6003
$ws.Range("M70").Value = -2992.5
$ws.Range("N70").Value = -6543
$ws.Range("H73").Value = 1635.6
$ws.Range("I73").Value = 1087.5
$ws.Range("J73").Value = 2001
$ws.Range("K73").Value = 3262.5
$ws.Range("L73").Value = 6003
$ws.Range("M73").Value = -2326.5
$ws.Range("N73").Value = -7875
$ws.Range("H125").Value = 6017.4546
$ws.Range("I125").Value = 6399.5
$ws.Range("K125").Value = 57595.5
$ws.Range("M125").Value = -55135.5
$ws.Range("H132").Value = 4551902.5
$ws.Range("I132").Value = 6623.684
$ws.Range("K132").Value = 19871.052
$ws.Range("M132").Value = -17341.052
$ws.Range("H137").Value = 6592.7144
$ws.Range("J137").Value = 2181.5557
$ws.Range("L137").Value = 6544.6671
$ws.Range("N137").Value = -11644.6671

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 4221.7744
$ws.Range("I32").Value = 3927.966
$ws.Range("K32").Value = 3927.966
$ws.Range("M32").Value = -3640.966
$ws.Range("H74").Value = 7464.7
$ws.Range("I74").Value = 7752.5293
$ws.Range("K74").Value = 7752.5293
$ws.Range("M74").Value = -6878.5293
$ws.Range("H77").Value = 7464.7
$ws.Range("I77").Value = 7752.5293
$ws.Range("K77").Value = 38762.6465
$ws.Range("M77").Value = -34394.6465

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H94").Value = 7235.0225
$ws.Range("I94").Value = 8893.362999999999
$ws.Range("K94").Value = 8893.362999999999
$ws.Range("M94").Value = -8442.362999999999
$ws.Range("H99").Value = 10616.774
$ws.Range("I99").Value = 11488.083
$ws.Range("K99").Value = 11488.083
$ws.Range("M99").Value = -9990.083000000001
$ws.Range("H105").Value = 71297.266
$ws.Range("I105").Value = 102885.9
$ws.Range("J105").Value = 8120
$ws.Range("K105").Value = 102885.9
$ws.Range("L105").Value = 8120
$ws.Range("M105").Value = -101138.9
$ws.Range("N105").Value = -11614
$ws.Range("H107").Value = 2522.4
$ws.Range("I107").Value = 2522.4
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2522.4
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -602.4000000000001
$ws.Range("N107").ClearContents()

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value = 1324.8334
$ws.Range("J22").Value = 1324.8334
$ws.Range("L22").Value = 1324.8334
$ws.Range("N22").Value = -2024.8334
$ws.Range("H31").Value = 10510.533
$ws.Range("I31").Value = 11050.615
$ws.Range("J31").Value = 7000
$ws.Range("K31").Value = 11050.615
$ws.Range("L31").Value = 7000
$ws.Range("M31").Value = -10755.615
$ws.Range("N31").Value = -7590
$ws.Range("H34").Value = 10510.533
$ws.Range("I34").Value = 11050.615
$ws.Range("J34").Value = 7000
$ws.Range("K34").Value = 11050.615
$ws.Range("L34").Value = 7000
$ws.Range("M34").Value = -10848.615
$ws.Range("N34").Value = -7404
$ws.Range("H122").Value = 12585
$ws.Range("I122").Value = 18650.428
$ws.Range("K122").Value = 55951.284
$ws.Range("M122").Value = -53501.284
$ws.Range("H132").Value = 1415.1072
$ws.Range("I132").Value = 1356.96
$ws.Range("K132").Value = 4070.88
$ws.Range("M132").Value = -1540.88
$ws.Range("H134").Value = 11202.363
$ws.Range("I134").Value = 14303.375
$ws.Range("K134").Value = 42910.125
$ws.Range("M134").Value = -40375.125
$ws.Range("H141").Value = 338343.28
$ws.Range("J141").Value = 411527.8
$ws.Range("L141").Value = 411527.8
$ws.Range("N141").Value = -421887.8

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H107").Value = 2018
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 2018
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 6054
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -9894
$ws.Range("H121").Value = 4669225
$ws.Range("I121").Value = 4222909
$ws.Range("K121").Value = 12668727
$ws.Range("M121").Value = -12667417
$ws.Range("J131").Value = 1983.7858
$ws.Range("L131").Value = 5951.357400000001
$ws.Range("N131").Value = -16031.3574
$ws.Range("H137").Value = 4324.643
$ws.Range("J137").Value = 13566
$ws.Range("L137").Value = 40698
$ws.Range("N137").Value = -50898

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H70").Value = 7656.6562
$ws.Range("I70").Value = 7018.619
$ws.Range("J70").Value = 8874.727999999999
$ws.Range("K70").Value = 7018.619
$ws.Range("L70").Value = 8874.727999999999
$ws.Range("M70").Value = -6748.619
$ws.Range("N70").Value = -9414.727999999999
$ws.Range("H73").Value = 7656.6562
$ws.Range("I73").Value = 7018.619
$ws.Range("J73").Value = 8874.727999999999
$ws.Range("K73").Value = 7018.619
$ws.Range("L73").Value = 8874.727999999999
$ws.Range("M73").Value = -6082.619
$ws.Range("N73").Value = -10746.728
$ws.Range("H102").Value = 5391.9165
$ws.Range("I102").Value = 5989.3706
$ws.Range("K102").Value = 5989.3706
$ws.Range("M102").Value = -4367.3706
$ws.Range("H107").Value = 614.94116
$ws.Range("I107").Value = 604.26666
$ws.Range("K107").Value = 604.26666
$ws.Range("M107").Value = 1315.73334
$ws.Range("H122").Value = 7495.4194
$ws.Range("I122").Value = 4133.1924
$ws.Range("J122").Value = 24979
$ws.Range("K122").Value = 12399.5772
$ws.Range("L122").Value = 74937
$ws.Range("M122").Value = -9949.5772
$ws.Range("N122").Value = -79837
$ws.Range("H126").Value = 10084.444
$ws.Range("I126").Value = 14178.8
$ws.Range("J126").Value = 7676
$ws.Range("K126").Value = 42536.39999999999
$ws.Range("L126").Value = 23028
$ws.Range("M126").Value = -40066.39999999999
$ws.Range("N126").Value = -27968
$ws.Range("H132").Value = 13014
$ws.Range("I132").Value = 14599.667
$ws.Range("K132").Value = 43799.001
$ws.Range("M132").Value = -41269.001

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H61").Value = 5178.8945
$ws.Range("I61").Value = 2069.5386
$ws.Range("J61").Value = 11915.833
$ws.Range("K61").Value = 2069.5386
$ws.Range("L61").Value = 11915.833
$ws.Range("M61").Value = -1867.5386
$ws.Range("N61").Value = -12319.833
$ws.Range("H74").Value = 36500
$ws.Range("J74").Value = 36500
$ws.Range("L74").Value = 36500
$ws.Range("N74").Value = -38496
$ws.Range("H77").Value = 36500
$ws.Range("J77").Value = 36500
$ws.Range("L77").Value = 109500
$ws.Range("N77").Value = -119484
$ws.Range("H113").Value = 5178.8945
$ws.Range("I113").Value = 2069.5386
$ws.Range("J113").Value = 11915.833
$ws.Range("K113").Value = 2069.5386
$ws.Range("L113").Value = 11915.833
$ws.Range("M113").Value = 100.4614000000001
$ws.Range("N113").Value = -16255.833
$ws.Range("H122").Value = 7197.2104
$ws.Range("I122").Value = 8041.9
$ws.Range("K122").Value = 24125.7
$ws.Range("M122").Value = -21675.7

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H126").Value = 33343.31
$ws.Range("I126").Value = 51796.375
$ws.Range("K126").Value = 155389.125
$ws.Range("M126").Value = -152919.125
